$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 3598093.74
$ws.Range("C7").Value = -19.01799511911351
$ws.Range("D7").Value = 3139
$ws.Range("E7").Value = 3139
$ws.Range("F7").Value = 1146.25477540618
$ws.Range("G7").Value = 22.18247056893228

$wb.Save()
